# Chapter2_Outline.docx edit script
# 1) Merge split runs "The Importance of Sea Ice " + "Thickness " + "initialization"
#    into a single run "The Importance of Sea Ice Thickness initialization".
# 2) Merge split runs "CanSIPSv1b " + "(Arlen)" into a single run "CanSIPSv1b (Arlen)".
# 3) Add 45 new character styles ListLabel46..ListLabel90 (same shape as the
#    existing ListLabel1..ListLabel45 styles: qFormat + rFonts cs="OpenSymbol").

$d = $word.ActiveDocument

# --- 1) & 2): collapse the split runs back into single runs ------------------
$null = $d.Content.Find.Execute(
    "The Importance of Sea Ice Thickness initialization",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The Importance of Sea Ice Thickness initialization", 2)

$null = $d.Content.Find.Execute(
    "CanSIPSv1b (Arlen)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CanSIPSv1b (Arlen)", 2)

# --- 3): mint the new ListLabel character styles ------------------------------
for ($i = 46; $i -le 90; $i++) {
    $style = $d.Styles.Add("ListLabel $i", 2)
    $style.QuickStyle = $true
    $style.Font.NameBi = "OpenSymbol"
}

Write-Output "done"
